$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Feuil1")

# --- Fill in the newly-graded rows of column A (style matches existing "1" cells, s=4) ---
$ws.Range("A4").Copy()
$greenRows = @(7,8,9,10,13,19,20,21)
foreach ($r in $greenRows) {
    $cell = $ws.Range("A$r")
    $cell.Value = 1
    $cell.PasteSpecial(-4122) | Out-Null
}

# --- A11 gets the new blue highlight style (new font + new fill) ---
$a11 = $ws.Range("A11")
$a11.Value = 1
$a11.Font.Color = 15773696
$a11.Interior.Color = 15773696

# --- Update the active selection on sheet1 ---
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B29").Select() | Out-Null

# --- Remove the now-unused empty sheets ---
$wb.Worksheets("Feuil2").Delete() | Out-Null
$wb.Worksheets("Feuil3").Delete() | Out-Null

# --- Rename the remaining sheet ---
$wb.Worksheets("Feuil1").Name = "Progression"
